$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '22.391.92'
$ws.Range("E2").Value = '  -0.35%  '
$ws.Range("D3").Value = '1.573.99'
$ws.Range("E3").Value = '  +0.03%  '
$ws.Range("E4").Value = '  -0.27%  '
$ws.Range("E5").Value = '  -0.24%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '291.07'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.30%  '
$ws.Range("E7").Value = '  +2.44%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '50.12'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.74%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3419'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.99%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.167'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.60%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07685'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.87%  '
$ws.Range("E12").Value = '  -0.20%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '21.39'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.88%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.996'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.48%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.936'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.39%  '
$ws.Range("E16").Value = '  +1.34%  '
$ws.Range("D17").Value = '1.575.26'
$ws.Range("E17").Value = '  +0.10%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '90.39'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.03%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06722'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.55%  '
$ws.Range("E20").Value = '  -0.19%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.79'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.21%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.247'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.28%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.5282'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.75%  '
$ws.Range("E24").Value = '  +1.62%  '
$ws.Range("D25").Value = '22.399.02'
$ws.Range("E25").Value = '  -0.36%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.390'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.16%  '
$ws.Range("E27").Value = '  -4.15%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.33'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.47%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '144.56'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.92%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.076'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.47%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '126.25'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.40%  '
$ws.Range("D32").Value = '1.748.40'
$ws.Range("E32").Value = '  -0.05%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.026'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +7.83%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.246'
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.023'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.25%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '10.08'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.43%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.08524'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.55%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02562'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.79%  '
$ws.Range("E39").Value = '  +1.88%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.06555'
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.520'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.06%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.295'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.42%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '11.65'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.14%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.6450'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.74%  '
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '14.17'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.34%  '
$ws.Range("B46").Value = 'Frax'
$ws.Range("C46").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.001'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.24%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.6023'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.28%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.780'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.06%  '
$ws.Range("E49").Value = '  +11.07%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.100'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.37%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '125.10'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.61%  '
